# The workbook tracks daily "Berenjena" (eggplant) price records. A new
# observation is being inserted as row 69 (pushing the existing rows 69-190
# down to 70-191), matching the structure/columns of its neighboring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 69; this shifts rows 69:190 down
# to 70:191 and carries the row's existing formatting (including the date
# style on column D) down with them.
$ws.Rows(69).Insert()

# Populate the newly inserted row 69 with the new record's data.
$ws.Range("A69").Value = 3
$ws.Range("B69").Value = "Femacal de La Calera"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 44533
$ws.Range("E69").Value = 5
$ws.Range("F69").Value = 100112001
$ws.Range("G69").Value = "Berenjena"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 120
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 8000
$ws.Range("M69").Value = 7729
$ws.Range("N69").Value = "$/caja 60 unidades"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 129
$ws.Range("Q69").Value = 60
$ws.Range("R69").Value = "Hortaliza"
